$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "A1" = 2204.2292172112393
    "B1" = 1386.0864356158613
    "C1" = 1434.8848532135837
    "A2" = 2227.1192368922416
    "B2" = 1482.0442356022584
    "C2" = 1329.6420231576942
    "A3" = 2343.9137069401977
    "B3" = 1532.1749178163932
    "C3" = 1629.9765938678086
    "A4" = 2320.9625339441868
    "B4" = 1791.2092367794739
    "C4" = 1735.3768306211518
    "A5" = 2422.6539933918984
    "B5" = 1663.7700490803886
    "C5" = 1643.2199901288009
    "A6" = 2360.7085097391273
    "B6" = 1775.9546966267449
    "C6" = 1784.4483163994612
    "A7" = 1992.9864161485818
    "B7" = 1556.0441516909277
    "C7" = 1533.3702716577748
    "A8" = 2135.7478628495105
    "B8" = 1655.7390017928153
    "C8" = 1507.358482697219
    "A9" = 2471.0153725375085
    "B9" = 1788.9612874546681
    "C9" = 1515.4135845837143
    "A10" = 2111.59440425195
    "B10" = 1363.7332583996726
    "C10" = 1304.7589921853487
    "A11" = 1970.0041216250017
    "B11" = 1414.7991874097338
    "C11" = 1308.366163587184
    "A12" = 2787.9537212083187
    "B12" = 2270.4390871133605
    "C12" = 2036.7338423286046
    "A13" = 2315.1853410978001
    "B13" = 1763.9087267538075
    "C13" = 1803.6382322348347
    "A14" = 2593.0481324833731
    "B14" = 1922.4492006468518
    "C14" = 1703.8047046772597
    "A15" = 2508.5108802462019
    "B15" = 2026.7229049020305
    "C15" = 1827.5083485611945
    "A16" = 2205.5222145157363
    "B16" = 1512.621046798469
    "C16" = 1374.8836059549424
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
